# Aggiornamento dati fino al 9 agosto 2021
# Appends rows 329-343 (dates 44403-44417) to the data table in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data: date serial, B (nuovi pos.), C (somma mobile 7gg.), D (somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(329, 44403, 0, 0, 0),
    @(330, 44404, 0, 0, 0),
    @(331, 44405, 0, 0, 0),
    @(332, 44406, 0, 0, 0),
    @(333, 44407, 2, 2, 37.87878787878788),
    @(334, 44408, 1, 3, 56.81818181818181),
    @(335, 44409, 0, 3, 56.81818181818181),
    @(336, 44410, 2, 5, 94.6969696969697),
    @(337, 44411, 1, 6, 113.6363636363636),
    @(338, 44412, 2, 8, 151.5151515151515),
    @(339, 44413, 1, 9, 170.4545454545454),
    @(340, 44414, 1, 8, 151.5151515151515),
    @(341, 44415, 0, 7, 132.5757575757576),
    @(342, 44416, 0, 7, 132.5757575757576),
    @(343, 44417, 3, 8, 151.5151515151515)
)

# Use the last existing row (328) as a formatting template for the new rows.
$templateRow = 328

foreach ($entry in $data) {
    $row = $entry[0]
    $dateSerial = $entry[1]
    $b = $entry[2]
    $c = $entry[3]
    $d = $entry[4]

    # Copy formatting from the template row into the new row first.
    $srcRange = $ws.Range("A$templateRow`:D$templateRow")
    $dstRange = $ws.Range("A$row`:D$row")
    $srcRange.Copy($dstRange)

    $ws.Cells.Item($row, 1).Value = $dateSerial
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
}
